$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds weekly price records for "Zapallo italiano" (Mapocho Venta
# Directa de Santiago). This edit refreshes the data: rows 2-12 (excluding the
# fixed descriptive columns A,B,C,E,F,G,H,I,R which stay constant) are
# re-populated with a re-shuffled set of weekly values - date (D), volume (J),
# min/max/avg price (K/L/M), unit of sale (N), origin (O), price per kg (P)
# and kg-or-units (Q).

# Column layout: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg,
# F Categoría ID, G Categoría, H Variedad, I Calidad, J Volumen,
# K Precio mínimo, L Precio máximo, M Precio promedio ponderado,
# N Unidad de comercialización, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificación.

$rows = @{
    2  = @{ D = 44277; J = 25; K = 10000; L = 10000; M = 10000; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';   P = 167; Q = 60 }
    3  = @{ D = 44179; J = 15; K = 7000;  L = 7000;  M = 7000;  N = '$/caja 60 unidades'; O = 'Provincia de Limarí';   P = 117; Q = 60 }
    4  = @{ D = 44405; J = 45; K = 9000;  L = 9000;  M = 9000;  N = '$/caja 50 unidades'; O = 'Provincia de Quillota'; P = 180; Q = 50 }
    5  = @{ D = 44200; J = 10; K = 9000;  L = 9000;  M = 9000;  N = '$/caja 60 unidades'; O = 'Provincia de Limarí';   P = 150; Q = 60 }
    7  = @{ D = 44315; J = 25; K = 10000; L = 10000; M = 10000; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';   P = 167; Q = 60 }
    8  = @{ D = 44291; J = 20; K = 9000;  L = 9000;  M = 9000;  N = '$/caja 60 unidades'; O = 'Provincia de Limarí';   P = 150; Q = 60 }
    9  = @{ D = 44186; J = 15; K = 7000;  L = 7000;  M = 7000;  N = '$/caja 60 unidades'; O = 'Provincia de Limarí';   P = 117; Q = 60 }
    11 = @{ D = 44333; J = 25; K = 10000; L = 11000; M = 10400; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';   P = 173; Q = 60 }
    12 = @{ D = 44312; J = 30; K = 10000; L = 10000; M = 10000; N = '$/caja 60 unidades'; O = 'Provincia de Limarí';   P = 167; Q = 60 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
}
